$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update total_risk (R) and total_risk_resp (S) values per newest airtoxics NATA data
$ws.Range("R2").Value = 60

$ws.Range("S3").Value = 0.5

$ws.Range("S4").Value = 0.35

$ws.Range("R5").Value = 28.3333333333333
$ws.Range("S5").Value = 0.316666666666667

$ws.Range("R6").Value = 30
$ws.Range("S6").Value = 0.35

$ws.Range("R7").Value = 60
$ws.Range("S7").Value = 0.425

$ws.Range("R8").Value = 60
$ws.Range("S8").Value = 0.5

$ws.Range("R9").Value = 30
$ws.Range("S9").Value = 0.3

$ws.Range("R11").Value = 16.6666666666667
